$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 12 currently holds the SUM totals (A12/B12 blank, C12 = SUM(C2:C11), D12 blank).
# We need to turn it into a new data row, and push the totals row down to row 13.

# 1. Copy the formatting of the current totals row (row 12) down onto row 13,
#    since that row keeps the same (border-only) formatting as the old totals row.
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D13").PasteSpecial($xlPasteFormats)

# 2. Copy the formatting of the last data row (row 11) onto the new data row 12.
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial($xlPasteFormats)

# 3. Fill in the new data row 12 values.
$ws.Cells.Item(12, 1).Value = 44319
$ws.Cells.Item(12, 2).Value = "oplossen problemen met Create, verschillende kleine aanpassingen aan design en views."
$ws.Cells.Item(12, 3).Value = 4
$ws.Cells.Item(12, 4).Value = ""

# 4. Move the totals formula down to row 13, extending the summed range.
$ws.Cells.Item(13, 1).Value = ""
$ws.Cells.Item(13, 2).Value = ""
$ws.Cells.Item(13, 3).Formula = "=SUM(C2:C12)"
$ws.Cells.Item(13, 4).Value = ""

$ws.Range("C13").Select()
